# Commit: "Fruta / hortaliza, semanal"
# Insert two new weekly price-report rows (Apio, Terminal Hortofrutícola
# Agro Chillán) ahead of the existing row 357, shifting the remaining
# data rows (357-404) down to (359-406).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 357-404 down by two rows, opening up two blank rows at 357:358.
$ws.Rows("357:358").Insert()

# --- New row 357: "Primera" quality -----------------------------------
$ws.Range("A357").Value = 7
$ws.Range("B357").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C357").Value = "Ñuble"
$ws.Range("D357").Value = 45131
$ws.Range("E357").Value = 16
$ws.Range("F357").Value = 100112017
$ws.Range("G357").Value = "Apio"
$ws.Range("H357").Value = "Americana (o)"
$ws.Range("I357").Value = "Primera"
$ws.Range("J357").Value = 200
$ws.Range("K357").Value = 7000
$ws.Range("L357").Value = 8000
$ws.Range("M357").Value = 7500
$ws.Range("N357").Value = "`$/docena de matas"
$ws.Range("O357").Value = "Provincia del Elquí"
$ws.Range("P357").Value = 1250
$ws.Range("Q357").Value = 6
$ws.Range("R357").Value = "Hortaliza"

# --- New row 358: "Segunda" quality -----------------------------------
$ws.Range("A358").Value = 7
$ws.Range("B358").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C358").Value = "Ñuble"
$ws.Range("D358").Value = 45131
$ws.Range("E358").Value = 16
$ws.Range("F358").Value = 100112017
$ws.Range("G358").Value = "Apio"
$ws.Range("H358").Value = "Americana (o)"
$ws.Range("I358").Value = "Segunda"
$ws.Range("J358").Value = 150
$ws.Range("K358").Value = 6000
$ws.Range("L358").Value = 6000
$ws.Range("M358").Value = 6000
$ws.Range("N358").Value = "`$/docena de matas"
$ws.Range("O358").Value = "Provincia del Elquí"
$ws.Range("P358").Value = 1000
$ws.Range("Q358").Value = 6
$ws.Range("R358").Value = "Hortaliza"

# Keep the date columns formatted like every other row in the column
# (style index 2 => "YYYY-MM-DD HH:MM:SS" custom date/time format).
$ws.Range("D357:D358").NumberFormat = $ws.Range("D356").NumberFormat

$wb.Save()
